# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 32 (Vega Central Mapocho de Santiago -
# Espárragos), pushing the existing rows 32:41 down to 33:42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 32:41 down to 33:42, leaving a fresh blank row 32.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new weekly observation.
$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(32, 3).Value = "Metropolitana"
$ws.Cells.Item(32, 4).Value = 44466
$ws.Cells.Item(32, 5).Value = 13
$ws.Cells.Item(32, 6).Value = 300000000
$ws.Cells.Item(32, 7).Value = "Espárragos"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Banquete"
$ws.Cells.Item(32, 10).Value = 110
$ws.Cells.Item(32, 11).Value = 2500
$ws.Cells.Item(32, 12).Value = 2500
$ws.Cells.Item(32, 13).Value = 2500
$ws.Cells.Item(32, 14).Value = "$/kilo"
$ws.Cells.Item(32, 15).Value = "Región Metropolitana"
$ws.Cells.Item(32, 16).Value = 2500
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Match the date cell formatting used by the rest of column D.
$ws.Cells.Item(32, 4).NumberFormat = $ws.Cells.Item(33, 4).NumberFormat
